$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Copy the highlighted "note" style (fillId 2, default number format —
#    the style already used by F2/F3/F4) onto the new note cells in
#    column F, then fill in their text. Re-copying per target cell avoids
#    flaky multi-area PasteSpecial behaviour.
# ---------------------------------------------------------------------
$noteCells = @("F29","F30","F45","F46","F47","F48","F49","F50","F54","F55","F56","F57","F58","F59")
foreach ($addr in $noteCells) {
    $ws.Range("F3").Copy()
    $ws.Range($addr).PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

$ws.Range("F29").Value = "Graphs do not show on this "
$ws.Range("F30").Value = "Graphs do not show on this "

$ws.Range("F45").Value = "Graphs do not show on this "
$ws.Range("F46").Value = "Graphs do not show on this "
$ws.Range("F47").Value = "Graphs do not show on this "

$ws.Range("F48").Value = "The plot is just flat and doesn’t look right"
$ws.Range("F49").Value = "The plot is just flat and doesn’t look right"
$ws.Range("F50").Value = "The plot is just flat and doesn’t look right"

$ws.Range("F54").Value = "The graph doesn’t show up"
$ws.Range("F55").Value = "The graph doesn’t show up"
$ws.Range("F56").Value = "The graph doesn’t show up"

$ws.Range("F57").Value = "The plot is just flat and doesn’t look right"
$ws.Range("F58").Value = "The plot is just flat and doesn’t look right"
$ws.Range("F59").Value = "The plot is just flat and doesn’t look right"

# ---------------------------------------------------------------------
# 2. G73 picks up the row's own highlighted style (fillId 4, like
#    A73:C73) rather than the column-F note style.
# ---------------------------------------------------------------------
$ws.Range("A73").Copy()
$ws.Range("G73").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("G73").Value = "Date in file says 7-13, while the file name says 7-8"

# ---------------------------------------------------------------------
# 3. New Time (column D) / ppm (column E) readings for rows that were
#    previously blank. D4 already has the correct "h:mm" + highlight
#    style (s=12) used throughout for Time readings in highlighted
#    blocks; E4 has the matching plain-highlight style (s=6) for ppm.
# ---------------------------------------------------------------------
$timeCells = @("D29","D30","D31","D89","D90","D91")
foreach ($addr in $timeCells) {
    $ws.Range("D4").Copy()
    $ws.Range($addr).PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

$ppmCells = @("E31","E89","E90","E91")
foreach ($addr in $ppmCells) {
    $ws.Range("E4").Copy()
    $ws.Range($addr).PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# D29 / D30 stay blank (formatted only); the rest get values.
$ws.Range("D31").Value = 0.4284722222222222
$ws.Range("E31").Value = 248.88140000000001

$ws.Range("D89").Value = 0.5180555555555556
$ws.Range("E89").Value = 1154.627

$ws.Range("D90").Value = 0.53194444444444444
$ws.Range("E90").Value = 1131.7460000000001

$ws.Range("D91").Value = 0.54236111111111118
$ws.Range("E91").Value = 1113.78

# ---------------------------------------------------------------------
# 4. Restore the view: scroll back to the top and select D92 (matches
#    the saved selection/scroll position in the edited workbook).
# ---------------------------------------------------------------------
$ws.Range("D92").Select() | Out-Null
